# Basic scoring functions for SSQNtotal finished; clear the stray
# leftover lookup cells (C2/D2/C3) that weren't supposed to be there yet,
# and leave the SSQStotal sheet's selection alone while switching focus
# back to SSQNtotal, which is where work continues.

$wb = $excel.ActiveWorkbook

$wsN = $wb.Worksheets.Item("SSQNtotal")
$wsS = $wb.Worksheets.Item("SSQStotal")

# Remove the not-yet-finished lookup values from SSQNtotal.
$wsN.Range("C2").ClearContents()
$wsN.Range("D2").ClearContents()
$wsN.Range("C3").ClearContents()

# Leave the cursor on SSQStotal where it was left off (B8) ...
$wsS.Activate()
$wsS.Range("B8").Select()

# ... then move on to SSQNtotal (now the active/selected sheet/tab),
# continuing work from cell C8.
$wsN.Activate()
$wsN.Range("C8").Select()
